# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newly calculated K values (column G, rows 2-36), replacing the stale
# Strike# derived numbers that were previously stored there.
$sVals = @(0, 4, 7, 6, 2, 8, 3, 1, 7, 2, 1, 6, 1, 4, 7, 5, 5, 9, 6, 7, 5, 7, 4, 5, 9, 4, 9, 2, 7, 10, 4, 3, 4, 5, 0)

$startRow = 2
for ($i = 0; $i -lt $sVals.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 7).Value = $sVals[$i]
}
